# Update "想去人数" (F column) counts on both the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 136
    3  = 1693
    5  = 28
    6  = 468
    7  = 154
    8  = 77
    9  = 614
    10 = 409
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
